$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for Price cells whose new values would otherwise
# be auto-parsed as numbers by Excel (losing significant trailing zeros
# or the exact textual representation used by the source site).
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply the updated Price / Volume(1h) text for each changed row.
$ws.Range("D2").Value = '30.101.70'
$ws.Range("E2").Value = '  +5.81%  '
$ws.Range("D3").Value = '1.922.54'
$ws.Range("E3").Value = '  +2.81%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.85%  '
$ws.Range("D5").Value = '331.68'
$ws.Range("E5").Value = '  +5.08%  '
$ws.Range("D6").Value = '0.9999'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").Value = '0.5251'
$ws.Range("E7").Value = '  +3.27%  '
$ws.Range("D8").Value = '0.4074'
$ws.Range("E8").Value = '  +4.53%  '
$ws.Range("D9").Value = '0.08546'
$ws.Range("E9").Value = '  +2.65%  '
$ws.Range("D10").Value = '43.44'
$ws.Range("E10").Value = '  +4.16%  '
$ws.Range("E11").Value = '  +2.59%  '
$ws.Range("D12").Value = '22.46'
$ws.Range("E12").Value = '  +10.17%  '
$ws.Range("D13").Value = '6.449'
$ws.Range("E13").Value = '  +3.66%  '
$ws.Range("D14").Value = '1.920.94'
$ws.Range("E14").Value = '  +2.51%  '
$ws.Range("D15").Value = '7.420'
$ws.Range("E15").Value = '  +2.21%  '
$ws.Range("D16").Value = '1.001'
$ws.Range("E16").Value = '  -0.91%  '
$ws.Range("D17").Value = '96.81'
$ws.Range("E17").Value = '  +6.20%  '
$ws.Range("E18").Value = '  +1.30%  '
$ws.Range("D19").Value = '0.06717'
$ws.Range("E19").Value = '  -0.06%  '
$ws.Range("D20").Value = '18.35'
$ws.Range("E20").Value = '  +3.81%  '
$ws.Range("D21").Value = '0.9993'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").Value = '6.068'
$ws.Range("E22").Value = '  +2.82%  '
$ws.Range("D23").Value = '30.113.37'
$ws.Range("E23").Value = '  +5.65%  '
$ws.Range("D24").Value = '11.32'
$ws.Range("E24").Value = '  +2.04%  '
$ws.Range("D25").Value = '2.229'
$ws.Range("E25").Value = '  +0.37%  '
$ws.Range("D26").Value = '2.142.93'
$ws.Range("E26").Value = '  +2.73%  '
$ws.Range("E27").Value = '  +3.01%  '
$ws.Range("D28").Value = '160.38'
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("E29").Value = '  +3.36%  '
$ws.Range("D30").Value = '129.58'
$ws.Range("E30").Value = '  +2.93%  '
$ws.Range("D31").Value = '1.085'
$ws.Range("E31").Value = '  +4.73%  '
$ws.Range("E32").Value = '  +1.80%  '
$ws.Range("D33").Value = '6.135'
$ws.Range("E33").Value = '  +6.77%  '
$ws.Range("E34").Value = '  +1.21%  '
$ws.Range("D35").Value = '0.02525'
$ws.Range("E35").Value = '  +3.02%  '
$ws.Range("D36").Value = '0.06620'
$ws.Range("E36").Value = '  +1.28%  '
$ws.Range("D37").Value = '0.2230'
$ws.Range("E37").Value = '  +3.53%  '
$ws.Range("D38").Value = '9.069'
$ws.Range("E38").Value = '  +2.96%  '
$ws.Range("E39").Value = '  +4.63%  '
$ws.Range("D40").Value = '5.227'
$ws.Range("E40").Value = '  +4.10%  '
$ws.Range("D41").Value = '0.6586'
$ws.Range("E41").Value = '  +3.43%  '
$ws.Range("E42").Value = '  +5.98%  '
$ws.Range("D43").Value = '1.244'
$ws.Range("D44").Value = '0.6211'
$ws.Range("E44").Value = '  +3.77%  '
$ws.Range("D45").Value = '13.39'
$ws.Range("E45").Value = '  +2.79%  '
$ws.Range("D46").Value = '3.792'
$ws.Range("E46").Value = '  +2.88%  '
$ws.Range("D47").Value = '2.096'
$ws.Range("E47").Value = '  +4.72%  '
$ws.Range("E48").Value = '  +2.90%  '
$ws.Range("D49").Value = '125.38'
$ws.Range("E49").Value = '  +3.01%  '
$ws.Range("D50").Value = '80.06'
$ws.Range("E50").Value = '  +5.07%  '
$ws.Range("D51").Value = '1.158'
$ws.Range("E51").Value = '  +1.06%  '
